$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A7").Value = "2025-07-14 18:18:00"
